$d = $word.ActiveDocument

# 1. Replace the intro GRP paragraph (paragraph 2)
$d.Paragraphs.Item(2).Range.Text = "The Gross Regional Product (GRP) serves as a vital economic indicator, reflecting the economic output of San Diego County. Over the period from 2019 to 2023, the county has demonstrated remarkable economic resilience and growth. In 2019, the total GRP was approximately `$244.28 billion. Despite the challenges posed by the COVID-19 pandemic, 2020 saw a slight increase to `$244.82 billion. The subsequent years marked a significant recovery and growth, with the GRP reaching `$268.87 billion in 2021, `$296.68 billion in 2022, and further climbing to `$308.71 billion in 2023. This upward trend underscores the county's robust economic health and adaptability."

# 2. Replace the "per capita" paragraph (paragraph 4) with the new "key industries" paragraph,
#    then insert a brand-new paragraph after it containing the new "per capita" text.
$d.Paragraphs.Item(4).Range.Text = "The growth in San Diego County's GRP is driven by several key industries. The government sector remains the largest contributor, with its share increasing from approximately `$45.19 billion in 2019 to `$52.92 billion in 2023. Manufacturing also plays a significant role, with its contribution rising from `$28.26 billion to `$31.67 billion over the same period. The professional, scientific, and technical services sector has shown remarkable growth, expanding from `$28.13 billion in 2019 to `$37.04 billion in 2023. Additionally, the health care and social assistance sector grew from `$15.38 billion to `$20.21 billion, while the finance and insurance sector increased from `$14.67 billion to `$19.50 billion. These sectors collectively highlight the diverse and dynamic nature of San Diego's economy."

$d.Paragraphs.Item(4).Range.InsertParagraphAfter()
$d.Paragraphs.Item(5).Range.Text = "The per capita GRP in San Diego County has also seen a steady increase, reflecting the distribution of economic benefits among the population. In 2019, the per capita GRP was approximately `$73,347, which rose slightly to `$74,278 in 2020. The following years saw more substantial increases, with per capita GRP reaching `$82,100 in 2021, `$90,557 in 2022, and `$94,916 in 2023. This growth in per capita GRP indicates an improvement in the economic well-being of the region's residents."

# 3. Table header cell: "Total GRP (billion $)" -> "Total GRP (Billion $)"
#    (re-fetch the table/cell each time since a structural edit can leave old handles,
#     and even the document's Paragraphs-by-index collection, stale)
$tbl = $d.Tables.Item(1)
$tbl.Cell(1, 2).Range.Text = "Total GRP (Billion `$)"

# 4. Table per-capita column values for 2020, 2021, 2022 ("-" -> actual numbers)
$tbl = $d.Tables.Item(1)
$tbl.Cell(3, 3).Range.Text = "74,278"
$tbl = $d.Tables.Item(1)
$tbl.Cell(4, 3).Range.Text = "82,100"
$tbl = $d.Tables.Item(1)
$tbl.Cell(5, 3).Range.Text = "90,557"

# 5 & 6. Caption text update and final concluding paragraph replacement.
# After the table-cell edits above, $d.Paragraphs (indexed collection) can go stale,
# so re-acquire paragraphs via $d.Content.Paragraphs, which stays accurate.
$d.Content.Paragraphs.Item(30).Range.Text = "San Diego County Total and Per Capita GRP (2019-2023)"
$d.Content.Paragraphs.Item(31).Range.Text = "In conclusion, San Diego County's economy has demonstrated resilience and consistent growth over the past five years. The steady increase in both total and per capita GRP highlights the region's economic vitality. Key sectors such as government, manufacturing, and professional services have been instrumental in driving this growth. The rising per capita GRP further suggests that the economic benefits are being effectively distributed among the population, enhancing the overall economic well-being of the county."
